$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.119.41'
$ws.Range("E2").Value = '  -2.36%  '
$ws.Range("D3").Value = '2.348.48'
$ws.Range("E3").Value = '  -3.19%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '85.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.97%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.528'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.03%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.483'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0809'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.14%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.93%  '
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '2.708.98'
$ws.Range("E13").Value = '  -3.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.43'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.81'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.93%  '
$ws.Range("D16").Value = '2.363.08'
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.761'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.59%  '
$ws.Range("D18").Value = '40.087.86'
$ws.Range("E18").Value = '  -2.25%  '
$ws.Range("D19").Value = '0.0₃0902'
$ws.Range("E19").Value = '  -2.53%  '
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '68.33'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.65'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.01%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.83'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -4.07%  '
$ws.Range("E29").Value = '  -3.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '34.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '153.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.50'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.33%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0719'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.60%  '
$ws.Range("E36").Value = '  -0.47%  '
$ws.Range("E37").Value = '  -5.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0992'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.68'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.73'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.67%  '
$ws.Range("D42").Value = '1.973.98'
$ws.Range("E42").Value = '  -1.28%  '
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0265'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.67%  '
$ws.Range("D48").Value = '2.567.57'
$ws.Range("E48").Value = '  -3.59%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '93.32'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.49'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.12%  '
$ws.Range("E51").Value = '  -3.22%  '
